# Fixed if multiple of the same setting
#
# The "Flight Mission Cycle" summary sheet lists each setting once with its
# "No. of cycles". This bug fix demonstrates/handles the case where the same
# setting (here "Typing") appears more than once: its count is corrected
# (2 -> 1 for the first occurrence) and a second row is added for the
# duplicate "Typing" entry with count 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Flight Mission Cycle")

# Make this the active sheet/tab (mirrors activeTab moving from "Typing" to
# "Flight Mission Cycle" and tabSelected moving sheets).
$ws.Activate() | Out-Null

# Correct the existing "Typing" row's cycle count.
$ws.Range("B2").Value = 1

# Add the duplicate "Typing" setting row.
$ws.Range("A4").Value = "Typing"
$ws.Range("B4").Value = 2

# Update the selection to match the new active cell on this sheet.
$ws.Range("E4").Select() | Out-Null
